$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.527.51"
$ws.Range("E2").Value = "  +1.75%  "

$ws.Range("D3").Value = "2.215.29"
$ws.Range("E3").Value = "  -0.26%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'240.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("E6").Value = "  -1.80%  "

$ws.Range("D7").Value = "'74.75"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.79%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("E9").Value = "  +1.46%  "

$ws.Range("D10").Value = "'41.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.33%  "

$ws.Range("D11").Value = "'0.0925"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.90%  "

$ws.Range("D12").Value = "'54.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.84%  "

$ws.Range("D13").Value = "'6.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.01%  "

$ws.Range("E14").Value = "  -1.73%  "

$ws.Range("D15").Value = "2.548.34"
$ws.Range("E15").Value = "  -0.16%  "

$ws.Range("D16").Value = "'14.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.98%  "

$ws.Range("D17").Value = "2.209.04"
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").Value = "'0.800"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.54%  "

$ws.Range("D19").Value = "42.386.05"
$ws.Range("E19").Value = "  +1.76%  "

$ws.Range("E20").Value = "  +0.71%  "

$ws.Range("D21").Value = "'70.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.62%  "

$ws.Range("D22").Value = "'5.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.76%  "

$ws.Range("D23").Value = "'9.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.90%  "

$ws.Range("D24").Value = "'229.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("D25").Value = "'2.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.29%  "

$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("E27").Value = "  -3.36%  "

$ws.Range("D28").Value = "'3.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.76%  "

$ws.Range("E29").Value = "  -2.00%  "

$ws.Range("D30").Value = "'172.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.51%  "

$ws.Range("D31").Value = "'36.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +21.10%  "

$ws.Range("E32").Value = "  -4.91%  "

$ws.Range("D33").Value = "'20.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.70%  "

$ws.Range("D34").Value = "'0.0791"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.21%  "

$ws.Range("D35").Value = "'5.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.76%  "

$ws.Range("E36").Value = "  -1.36%  "

$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("D38").Value = "'4.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.83%  "

$ws.Range("E39").Value = "  +6.48%  "

$ws.Range("D40").Value = "'12.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.54%  "

$ws.Range("E41").Value = "  +1.38%  "

$ws.Range("E42").Value = "  -2.04%  "

$ws.Range("D43").Value = "'60.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.77%  "

$ws.Range("E44").Value = "  +0.80%  "

$ws.Range("E45").Value = "  -0.72%  "

$ws.Range("D46").Value = "'0.0988"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("D47").Value = "'99.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.72%  "

$ws.Range("E48").Value = "  -1.32%  "

$ws.Range("E49").Value = "  -1.60%  "

$ws.Range("E50").Value = "  -2.50%  "

$ws.Range("D51").Value = "'0.424"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +16.38%  "
